$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $pattern) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $pp = $doc.Paragraphs.Item($i)
        if ($pp.Range.Text -like $pattern) {
            return $pp
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1. Move the "_GoBack" bookmark from the end of the document (after the
#    last paragraph's text) to the start of the "Objetivos Funcionales"
#    heading paragraph (right after its pPr, before the run).
# ---------------------------------------------------------------------------
$bmOld = $d.Bookmarks.Item("_GoBack")
$endPos = $bmOld.End
$bmOld.Delete()

$headingPara = Find-ParagraphByText $d "Objetivos Funcionales*"
$startPos = $headingPara.Range.Start
$d.Bookmarks.Add("_GoBack", $d.Range($startPos, $endPos)) | Out-Null

# ---------------------------------------------------------------------------
# 2. Insert a new paragraph "Conectar con social media network" right
#    before "Implementar la gestión de perfiles." — "network" is wrapped in
#    spell-check (proofErr) markers, like the existing "login"/"multi-auth"
#    runs elsewhere in the document.
# ---------------------------------------------------------------------------
$perfilesPara = Find-ParagraphByText $d "*gestión de perfiles.*"
$perfilesPara.Range.InsertParagraphBefore()
# After InsertParagraphBefore, the new empty paragraph immediately precedes
# $perfilesPara (re-fetch it since the paragraph collection shifted).
$perfilesPara = Find-ParagraphByText $d "*gestión de perfiles.*"
$socialPara = $perfilesPara.Previous()

$socialXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Conectar con social media </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>network</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:sectPr/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$socialPara.Range.InsertXML($socialXml) | Out-Null

# ---------------------------------------------------------------------------
# 3. Insert a new paragraph "Implementar gestión de calificaciones de
#    usuarios" right before "Implementar la gestión de artículos
#    comprados."
# ---------------------------------------------------------------------------
$compradosPara = Find-ParagraphByText $d "*artículos comprados.*"
$compradosPara.Range.InsertBefore("Implementar gestión de calificaciones de usuarios`r")

# ---------------------------------------------------------------------------
# 4. Move the <w:lastRenderedPageBreak/> marker from the "adquisición"
#    paragraph to the "comprados" paragraph (re-fetch, text shifted nothing
#    but indices moved).
# ---------------------------------------------------------------------------
$adqPara = Find-ParagraphByText $d "*adquisición de artículos.*"
$compradosPara = Find-ParagraphByText $d "*artículos comprados.*"

$adqXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00377A70" w:rsidRDefault="00377A70"><w:r><w:t>Implementar la adquisición de artículos.</w:t></w:r></w:p><w:sectPr/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$adqPara.Range.InsertXML($adqXml) | Out-Null

$compradosXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00377A70" w:rsidRDefault="00377A70"><w:r><w:lastRenderedPageBreak/><w:t>Implementar la gestión de artículos comprados.</w:t></w:r></w:p><w:sectPr/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$compradosPara.Range.InsertXML($compradosXml) | Out-Null

Write-Output "Edit complete"
